# Update statistical description values (Mean, STD, quartiles, etc.)
# recomputed after data querying / ENV fLDPLM fitting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 517.416566736907
$ws.Range("D2").Value = 118.2430090347994
$ws.Range("F2").Value = 441
$ws.Range("G2").Value = 473
$ws.Range("H2").Value = 553
$ws.Range("C3").Value = 37.08015121039799
$ws.Range("D3").Value = 6.48922326979137
$ws.Range("F3").Value = 32.36
$ws.Range("G3").Value = 37.22
$ws.Range("H3").Value = 41.42
$ws.Range("I3").Value = 62.81
$ws.Range("C4").Value = 2.103507656356632
$ws.Range("D4").Value = 2.580611799609351
$ws.Range("F4").Value = 0.7
$ws.Range("G4").Value = 1.36
$ws.Range("H4").Value = 2.58
$ws.Range("C5").Value = 322.5873181168558
$ws.Range("D5").Value = 8.773077343613297
$ws.Range("F5").Value = 317.43
$ws.Range("G5").Value = 323.09
$ws.Range("H5").Value = 328.75
$ws.Range("C6").Value = 23.68275723078164
$ws.Range("D6").Value = 3.79232396646391
$ws.Range("F6").Value = 20.91
$ws.Range("G6").Value = 23.16
$ws.Range("H6").Value = 26.38
$ws.Range("C7").Value = -75.55727063374835
$ws.Range("D7").Value = 22.23745506296449
$ws.Range("F7").Value = -91
$ws.Range("C8").Value = 7.828773339471671
$ws.Range("D8").Value = 6.672697927244556
$ws.Range("C9").Value = 9.162045139899666
$ws.Range("D9").Value = 1.6243388762162
$ws.Range("C10").Value = 867.8246188992435
$ws.Range("D10").Value = 0.4612275166605061
$ws.Range("C11").Value = 0.4917044040828045
$ws.Range("D11").Value = 0.5482100045929824
$ws.Range("C12").Value = 22.74950993487495
$ws.Range("D12").Value = 12.29619773727692
$ws.Range("C13").Value = 0.6722120196034036
$ws.Range("D13").Value = 0.7501695231887071
$ws.Range("C14").Value = 1.830738753459606
$ws.Range("D14").Value = 1.668290349302445
$ws.Range("C15").Value = 92.81727063374788
$ws.Range("D15").Value = 22.2374550630588
$ws.Range("H15").Value = 108.26
$ws.Range("C16").Value = -84.9017073378634
$ws.Range("D16").Value = 20.02160869245486
$ws.Range("F16").Value = -100.9574620641016
$ws.Range("G16").Value = -82.46683163887967
$ws.Range("H16").Value = -69.0778545523916
$ws.Range("C17").Value = -77.07293399839182
$ws.Range("D17").Value = 24.544925109817
$ws.Range("F17").Value = -91.6389203414338
$ws.Range("G17").Value = -72.2376019773414
$ws.Range("H17").Value = -58.3175485570292
